$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C3").Value = "Gabriel"
$ws.Range("F3").Value = "Em andamento"

$ws.Range("C7").Value = "Gilvanildo"
$ws.Range("F7").Value = "Em andamento"

$ws.Range("A1:F1").Select()
